$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P&M Schedule")

# New rows 12-20: Table column = "AddDataPoint(P&M)", Column = <new field name>
$rows = @(
    "ProjectID",
    "ForecastCompletionYear",
    "ProjectType",
    "ProjectStage",
    "Criticality",
    "DelayInSchedule",
    "CostOverrun",
    "Priority",
    "ProjectValue"
)

$r = 12
foreach ($name in $rows) {
    $ws.Range("A$r").Value = "AddDataPoint(P&M)"
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $name
    $r = $r + 1
}

$ws.Range("B24").Select()
$ws.Activate()
